$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.357000000000001
$ws.Range("D4").Value = -7.961
$ws.Range("D7").Value = -8.120999999999999
$ws.Range("A9").Value = -21.723
$ws.Range("B9").Value = 5.554
$ws.Range("C9").Value = -11.319
$ws.Range("B11").Value = 5.909000000000001
$ws.Range("D11").Value = -7.562
$ws.Range("D15").Value = -8.024000000000001
$ws.Range("A18").Value = -21.632
$ws.Range("A20").Value = -20.338
$ws.Range("B23").Value = 7.542
$ws.Range("B24").Value = 6.359
$ws.Range("B26").Value = 6.308999999999999
$ws.Range("A27").Value = -21.886
$ws.Range("C27").Value = -13.145
$ws.Range("C29").Value = -12.248
$ws.Range("D30").Value = -7.450999999999999
$ws.Range("C32").Value = -12.079
$ws.Range("B34").Value = 7.129
$ws.Range("A35").Value = -21.689
$ws.Range("B35").Value = 5.577
$ws.Range("C37").Value = -12.266
$ws.Range("C38").Value = -12.209
$ws.Range("D38").Value = -7.975
$ws.Range("D39").Value = -7.696
$ws.Range("C41").Value = -12.375
$ws.Range("D43").Value = -7.441
$ws.Range("C45").Value = -13.297
$ws.Range("D47").Value = -7.515000000000001
$ws.Range("B48").Value = 6.329000000000001
$ws.Range("B49").Value = 6.317
$ws.Range("C51").Value = -11.356
$ws.Range("B52").Value = 5.766
$ws.Range("C57").Value = -13.843
$ws.Range("C64").Value = -11.016
$ws.Range("B66").Value = 5.652
$ws.Range("B67").Value = 5.718999999999999
$ws.Range("A69").Value = -21.408
$ws.Range("D75").Value = -7.447
$ws.Range("A76").Value = -20.564
$ws.Range("A78").Value = -20.996
$ws.Range("B78").Value = 6.561
$ws.Range("B80").Value = 7.122
$ws.Range("A82").Value = -21.408
$ws.Range("C82").Value = -11.8
$ws.Range("A83").Value = -21.852
$ws.Range("D91").Value = -7.669999999999999
$ws.Range("D92").Value = -7.003
$ws.Range("A93").Value = -21.377
$ws.Range("C93").Value = -11.25
$ws.Range("D95").Value = -7.358
$ws.Range("D96").Value = -7.211
$ws.Range("B99").Value = 5.645
$ws.Range("C102").Value = -12.788
$ws.Range("B104").Value = 7.550999999999999
$ws.Range("C105").Value = -12.39
$ws.Range("D105").Value = -7.495
